# localization-status.xlsx — "Generate Report for Archive"
#
# 1. Flip the in-progress rows from "Ready for handoff" to "In Translation"
#    on every sheet that carries that status (Overview + each locale sheet).
# 2. Tighten the now-shorter status column so it isn't over-wide.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newWidth  = 12.5   # renders as the narrower "status" column width

# --- Overview sheet: status appears in columns E (zh-cn) and F (de-de) ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($addr in @("E2","F2","E3","F3")) {
    $cell = $overview.Range($addr)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# --- Locale detail sheets (zh-cn, de-de): status is column C ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in @("C2","C3")) {
        $cell = $ws.Range($addr)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    $ws.Columns.Item(3).ColumnWidth = $newWidth
}
